# This script swaps the data content between specific pairs of rows in the
# "Artfynd" worksheet, matching the target diff. Each pair of rows keeps its
# shared/common column values (location, date, observer, etc.) in place but
# exchanges the record-specific values (id, taxon info, coordinates, and any
# extra per-record columns) between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($Row1, $Row2, $Columns)

    foreach ($col in $Columns) {
        $addr1 = "$col$Row1"
        $addr2 = "$col$Row2"

        $range1 = $ws.Range($addr1)
        $range2 = $ws.Range($addr2)

        $val1 = $range1.Value2
        $val2 = $range2.Value2

        # Treat "no value" the same whether it comes back as $null or an
        # empty string so we correctly clear cells that should end up blank.
        $has1 = -not ([string]::IsNullOrEmpty($val1))
        $has2 = -not ([string]::IsNullOrEmpty($val2))

        if ($has2) {
            $range1.Value = $val2
        } else {
            $range1.ClearContents()
        }

        if ($has1) {
            $range2.Value = $val1
        } else {
            $range2.ClearContents()
        }
    }
}

# Columns that differ between row 14 and row 15 (and likewise for the other
# simple pairs 17/19 and 30/31): identifying id, taxon data, and coordinates.
$commonColumns = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

Swap-RowData 14 15 $commonColumns
Swap-RowData 17 19 $commonColumns
Swap-RowData 30 31 $commonColumns

# Row 61 / 62 pair also carries extra columns (Ålder-Stadium, Kön, Aktivitet,
# Metod, Publik kommentar) that belong to only one of the two records.
$extendedColumns = @("A", "B", "D", "E", "F", "G", "H", "K", "L", "M", "N", "Q", "R", "AC")
Swap-RowData 61 62 $extendedColumns
